$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 9: C9, C10, C11, C12 -> C9-C12 (designator notation only, qty stays 4) ---
$ws.Range("A9").Value = "C9-C12"

# --- Row 8: D1, D2, D3 -> D1, D2, D3, D4 (now includes a 4th diode), qty 3 -> 4 ---
$ws.Range("A8").Value = "D1, D2, D3, D4"
$ws.Range("B8").Value = 4

# --- Row 18: R11 -> R11, R14 (added resistor), qty 1 -> 2 ---
$ws.Range("A18").Value = "R11, R14"
$ws.Range("B18").Value = 2

# --- New row 22: TP1, TP2 test points ---
$ws.Range("A22").Value = "TP1, TP2"
$ws.Range("B22").Value = 2
$ws.Range("D22").Value = "36-5000-ND"
$ws.Range("E22").Value = "DK"

# Add the hyperlink for the new part number, then restore the usual
# bordered/hyperlink cell formatting (matching the other Part Number cells)
# by copying formats from an existing, identically-styled Part Number cell.
$ws.Hyperlinks.Add($ws.Range("D22"), "https://www.digikey.com/en/products/detail/keystone-electronics/5000/316860") | Out-Null
$ws.Range("D3").Copy() | Out-Null
$ws.Range("D22").PasteSpecial(-4122) | Out-Null
$ws.Range("D22").Value = "36-5000-ND"

# --- Tidy up formatting left over from the removed/unused style (D21 had a
# stray non-bordered hyperlink style previously; make it consistent with the
# rest of column D) ---
$ws.Range("D3").Copy() | Out-Null
$ws.Range("D21").PasteSpecial(-4122) | Out-Null
$ws.Range("D21").Value = "CSTNE12M0GH5L000R0"

# --- Update the active selection to reflect where editing left off ---
$ws.Range("B23").Select() | Out-Null
